$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells in column D may look numeric, so force them to stay text
# using a Text number format then resetting the style to avoid altering
# the visual appearance, mirroring the original inline-string cells.
$textCells = @(
    @('D2', '26.918.66'),
    @('D3', '1.871.40'),
    @('D4', '0.9994'),
    @('D5', '304.77'),
    @('D6', '0.9991'),
    @('D7', '0.5088'),
    @('D8', '0.3668'),
    @('D9', '0.07182'),
    @('D11', '20.64'),
    @('D12', '0.07493'),
    @('D13', '1.878.71'),
    @('D14', '94.59'),
    @('D16', '0.9996'),
    @('D17', '0.000008511'),
    @('D18', '14.15'),
    @('D19', '0.9992'),
    @('D20', '26.959.27'),
    @('D22', '2.115.63'),
    @('D24', '6.391'),
    @('D26', '1.775'),
    @('D28', '2.083'),
    @('D29', '113.56'),
    @('D30', '4.691'),
    @('D31', '4.713'),
    @('D33', '0.05055'),
    @('D34', '0.7507'),
    @('D36', '1.153'),
    @('D37', '3.213'),
    @('D38', '0.5649'),
    @('D39', '2.523'),
    @('D42', '6.615'),
    @('D43', '115.47'),
    @('D44', '8.545'),
    @('D46', '0.4781'),
    @('D47', '0.9989'),
    @('D48', '10.08')
)

foreach ($pair in $textCells) {
    $cell = $ws.Range($pair[0])
    $cell.NumberFormat = "@"
    $cell.Value = $pair[1]
    $cell.Style = "Normal"
}

# Cells in columns B, C and E are never numeric-looking (URLs, coin
# names or whitespace-padded percentages), so a plain value assignment
# is sufficient and keeps the original default styling untouched.
$plainCells = @(
    @('E2', '  -0.61%  '),
    @('E3', '  +0.39%  '),
    @('E4', '  -0.18%  '),
    @('E5', '  -0.44%  '),
    @('E6', '  -0.15%  '),
    @('E7', '  -1.22%  '),
    @('E8', '  -2.40%  '),
    @('E9', '  +0.35%  '),
    @('E10', '  +0.23%  '),
    @('E11', '  -0.20%  '),
    @('B12', 'TRON'),
    @('C12', 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'),
    @('E12', '  -0.85%  '),
    @('B13', 'WrappedEther'),
    @('C13', 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'),
    @('E13', '  +0.88%  '),
    @('E14', '  +5.53%  '),
    @('E15', '  -1.51%  '),
    @('E16', '  -0.24%  '),
    @('E17', '  +0.44%  '),
    @('E18', '  +0.69%  '),
    @('E19', '  -0.10%  '),
    @('E20', '  -0.59%  '),
    @('E21', '  +0.08%  '),
    @('E22', '  +1.87%  '),
    @('E23', '  -1.12%  '),
    @('E24', '  -0.77%  '),
    @('E25', '  +1.64%  '),
    @('E26', '  -3.48%  '),
    @('E27', '  -0.40%  '),
    @('E28', '  -0.16%  '),
    @('E29', '  +0.67%  '),
    @('E30', '  +0.67%  '),
    @('E31', '  +1.22%  '),
    @('E32', '  +0.09%  '),
    @('E33', '  -0.88%  '),
    @('E34', '  +3.85%  '),
    @('E35', '  -2.95%  '),
    @('E36', '  -0.04%  '),
    @('E37', '  +4.14%  '),
    @('B38', 'TheSandbox'),
    @('C38', 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'),
    @('E38', '  +7.13%  '),
    @('B39', 'RenderToken'),
    @('C39', 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'),
    @('E39', '  +1.51%  '),
    @('E40', '  -2.24%  '),
    @('E41', '  -0.33%  '),
    @('E42', '  +2.32%  '),
    @('E43', '  -0.50%  '),
    @('E44', '  +3.30%  '),
    @('E45', '  +1.28%  '),
    @('E46', '  +3.62%  '),
    @('E47', '  -0.16%  '),
    @('E48', '  +1.41%  '),
    @('E49', '  -0.17%  '),
    @('E50', '  +1.29%  '),
    @('E51', '  -0.28%  ')
)

foreach ($pair in $plainCells) {
    $ws.Range($pair[0]).Value = $pair[1]
}
